$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated save_data: column G ("K") now reflects strikeouts (K) instead of
# the old Strike# metric. Write the recalculated s_vals for rows 2-34.
$kValues = @(0,2,2,2,1,1,0,2,3,1,2,0,2,0,3,1,0,1,1,2,2,1,2,1,0,1,1,0,2,2,0,1,2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}

$wb.Save()
